$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("InventoryQuantity"))
$ws.Name = "Prepayment"
$ws.Range("A1").Value = "Division"
Write-Output "done"
